$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 14.53326577974471

$ws.Range("B3").Value = 3.230985683306322
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 39.7764191927396

$ws.Range("B4").Value = 0.6753301551942219
$ws.Range("C4").Value = 10.29869402782916
$ws.Range("D4").Value = 26.21740644021617
$ws.Range("E4").Value = 8.660232485948974
$ws.Range("G4").Value = 45.85166310918853
